$wb = $excel.ActiveWorkbook

# ---- Overview ----
$ws1 = $wb.Worksheets.Item("Overview")

# Cell values
$ws1.Range("A2").Value = "fffff416882d-0faf-47cd-8774-c7dfb82d5749.md"
$ws1.Range("B2").Value = "Handed back: in sync with en-US"
$ws1.Range("C2").Value = "Handed back: in sync with en-US"
$ws1.Range("D2").Value = "2016-03-21 08:56:49"
$ws1.Range("A3").Value = "ffffff4298431a-edd8-4e9e-8464-d1e2453c7eec.md"
$ws1.Range("B3").Value = "Handed back: in sync with en-US"
$ws1.Range("C3").Value = "Handed back: in sync with en-US"
$ws1.Range("D3").Value = "2016-03-21 08:56:49"
$ws1.Range("A4").Value = "3343e4d7-024a-4224-a054-9d6db3bbcd44.md"
$ws1.Range("B4").Value = "Ready for handoff"
$ws1.Range("C4").Value = "Ready for handoff"
$ws1.Range("D4").Value = "2016-03-21 09:00:11"

# Hyperlinks: clear existing, re-add with original targets + refreshed display text
$ws1.Cells.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/98d1c6365ca3582fab6e8cfc607746a763847726/e2e/3343e4d7-024a-4224-a054-9d6db3bbcd44.md", "", "", "fffff416882d-0faf-47cd-8774-c7dfb82d5749.md")
$ws1.Hyperlinks.Add($ws1.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/98d1c6365ca3582fab6e8cfc607746a763847726/e2e/fffff416882d-0faf-47cd-8774-c7dfb82d5749.md", "", "", "ffffff4298431a-edd8-4e9e-8464-d1e2453c7eec.md")
$ws1.Hyperlinks.Add($ws1.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/98d1c6365ca3582fab6e8cfc607746a763847726/e2e/ffffff4298431a-edd8-4e9e-8464-d1e2453c7eec.md", "", "", "3343e4d7-024a-4224-a054-9d6db3bbcd44.md")

# ---- zh-cn ----
$ws2 = $wb.Worksheets.Item("zh-cn")

# Cell values
$ws2.Range("A2").Value = "fffff416882d-0faf-47cd-8774-c7dfb82d5749.md"
$ws2.Range("B2").Value = ".md"
$ws2.Range("C2").Value = "Handed back: in sync with en-US"
$ws2.Range("D2").Value = "227dfc44-15b6-42d6-a696-ae27c012d273.49a318ecf0566613ef68ad4b941dea8039fc0f09.zh-cn.xlf"
$ws2.Range("E2").Value = "2016-03-21 08:56:45"
$ws2.Range("F2").Value = "227dfc44-15b6-42d6-a696-ae27c012d273.md"
$ws2.Range("G2").Value = "227dfc44-15b6-42d6-a696-ae27c012d273.49a318ecf0566613ef68ad4b941dea8039fc0f09.zh-cn.xlf"
$ws2.Range("H2").Value = "2016-03-21 08:57:12"
$ws2.Range("J2").Value = "Include"
$ws2.Range("A3").Value = "ffffff4298431a-edd8-4e9e-8464-d1e2453c7eec.md"
$ws2.Range("B3").Value = ".md"
$ws2.Range("C3").Value = "Handed back: in sync with en-US"
$ws2.Range("D3").Value = "227dfc44-15b6-42d6-a696-ae27c012d273.49a318ecf0566613ef68ad4b941dea8039fc0f09.zh-cn.xlf"
$ws2.Range("E3").Value = "2016-03-21 08:56:45"
$ws2.Range("F3").Value = "227dfc44-15b6-42d6-a696-ae27c012d273.md"
$ws2.Range("G3").Value = "227dfc44-15b6-42d6-a696-ae27c012d273.49a318ecf0566613ef68ad4b941dea8039fc0f09.zh-cn.xlf"
$ws2.Range("H3").Value = "2016-03-21 08:57:12"
$ws2.Range("J3").Value = "Include"
$ws2.Range("A4").Value = "3343e4d7-024a-4224-a054-9d6db3bbcd44.md"
$ws2.Range("B4").Value = ".md"
$ws2.Range("C4").Value = "Ready for handoff"
$ws2.Range("D4").Value = "3343e4d7-024a-4224-a054-9d6db3bbcd44.fbdc8a0c47139da47674a6de4d6b421568f325ef.zh-cn.xlf"
$ws2.Range("E4").Value = "2016-03-21 09:00:05"
$ws2.Range("F4").Value = "3343e4d7-024a-4224-a054-9d6db3bbcd44.md"
$ws2.Range("G4").Value = "3343e4d7-024a-4224-a054-9d6db3bbcd44.fbdc8a0c47139da47674a6de4d6b421568f325ef.zh-cn.xlf"
$ws2.Range("H4").Value = "2016-03-21 08:59:30"
$ws2.Range("J4").Value = "Include"

# Hyperlinks: clear existing, re-add with original targets + refreshed display text
$ws2.Cells.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/98d1c6365ca3582fab6e8cfc607746a763847726/e2e/3343e4d7-024a-4224-a054-9d6db3bbcd44.md", "", "", "fffff416882d-0faf-47cd-8774-c7dfb82d5749.md")
$ws2.Hyperlinks.Add($ws2.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/fc3dbc4d1941dfd52c7540f5e2661b2aff760e8b/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/3343e4d7-024a-4224-a054-9d6db3bbcd44.fbdc8a0c47139da47674a6de4d6b421568f325ef.zh-cn.xlf", "", "", "227dfc44-15b6-42d6-a696-ae27c012d273.49a318ecf0566613ef68ad4b941dea8039fc0f09.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/1aa71b564f76bcc5f010e635e267ed2758a27d39/e2e/3343e4d7-024a-4224-a054-9d6db3bbcd44.md", "", "", "227dfc44-15b6-42d6-a696-ae27c012d273.md")
$ws2.Hyperlinks.Add($ws2.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/6e5735860ebaf3a8ea29392f583b973b1695c38c/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/3343e4d7-024a-4224-a054-9d6db3bbcd44.fbdc8a0c47139da47674a6de4d6b421568f325ef.zh-cn.xlf", "", "", "227dfc44-15b6-42d6-a696-ae27c012d273.49a318ecf0566613ef68ad4b941dea8039fc0f09.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/98d1c6365ca3582fab6e8cfc607746a763847726/e2e/fffff416882d-0faf-47cd-8774-c7dfb82d5749.md", "", "", "ffffff4298431a-edd8-4e9e-8464-d1e2453c7eec.md")
$ws2.Hyperlinks.Add($ws2.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5dd169b752834a2f3d6b9be2527ab9ae06400fc5/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/227dfc44-15b6-42d6-a696-ae27c012d273.49a318ecf0566613ef68ad4b941dea8039fc0f09.zh-cn.xlf", "", "", "227dfc44-15b6-42d6-a696-ae27c012d273.49a318ecf0566613ef68ad4b941dea8039fc0f09.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("F3"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/d098a511b843d27c812cf5f1a81b32eaec13e35b/e2e/227dfc44-15b6-42d6-a696-ae27c012d273.md", "", "", "227dfc44-15b6-42d6-a696-ae27c012d273.md")
$ws2.Hyperlinks.Add($ws2.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/db3c37e8ac2df5c155bfeeec5c1dc0369ed7733e/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/227dfc44-15b6-42d6-a696-ae27c012d273.49a318ecf0566613ef68ad4b941dea8039fc0f09.zh-cn.xlf", "", "", "227dfc44-15b6-42d6-a696-ae27c012d273.49a318ecf0566613ef68ad4b941dea8039fc0f09.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/98d1c6365ca3582fab6e8cfc607746a763847726/e2e/ffffff4298431a-edd8-4e9e-8464-d1e2453c7eec.md", "", "", "3343e4d7-024a-4224-a054-9d6db3bbcd44.md")
$ws2.Hyperlinks.Add($ws2.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5dd169b752834a2f3d6b9be2527ab9ae06400fc5/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/227dfc44-15b6-42d6-a696-ae27c012d273.49a318ecf0566613ef68ad4b941dea8039fc0f09.zh-cn.xlf", "", "", "3343e4d7-024a-4224-a054-9d6db3bbcd44.fbdc8a0c47139da47674a6de4d6b421568f325ef.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("F4"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/d098a511b843d27c812cf5f1a81b32eaec13e35b/e2e/227dfc44-15b6-42d6-a696-ae27c012d273.md", "", "", "3343e4d7-024a-4224-a054-9d6db3bbcd44.md")
$ws2.Hyperlinks.Add($ws2.Range("G4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/db3c37e8ac2df5c155bfeeec5c1dc0369ed7733e/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/227dfc44-15b6-42d6-a696-ae27c012d273.49a318ecf0566613ef68ad4b941dea8039fc0f09.zh-cn.xlf", "", "", "3343e4d7-024a-4224-a054-9d6db3bbcd44.fbdc8a0c47139da47674a6de4d6b421568f325ef.zh-cn.xlf")

# ---- de-de ----
$ws3 = $wb.Worksheets.Item("de-de")

# Cell values
$ws3.Range("A2").Value = "fffff416882d-0faf-47cd-8774-c7dfb82d5749.md"
$ws3.Range("B2").Value = ".md"
$ws3.Range("C2").Value = "Handed back: in sync with en-US"
$ws3.Range("D2").Value = "227dfc44-15b6-42d6-a696-ae27c012d273.49a318ecf0566613ef68ad4b941dea8039fc0f09.de-de.xlf"
$ws3.Range("E2").Value = "2016-03-21 08:56:49"
$ws3.Range("F2").Value = "227dfc44-15b6-42d6-a696-ae27c012d273.md"
$ws3.Range("G2").Value = "227dfc44-15b6-42d6-a696-ae27c012d273.49a318ecf0566613ef68ad4b941dea8039fc0f09.de-de.xlf"
$ws3.Range("H2").Value = "2016-03-21 08:57:19"
$ws3.Range("J2").Value = "Include"
$ws3.Range("A3").Value = "ffffff4298431a-edd8-4e9e-8464-d1e2453c7eec.md"
$ws3.Range("B3").Value = ".md"
$ws3.Range("C3").Value = "Handed back: in sync with en-US"
$ws3.Range("D3").Value = "227dfc44-15b6-42d6-a696-ae27c012d273.49a318ecf0566613ef68ad4b941dea8039fc0f09.de-de.xlf"
$ws3.Range("E3").Value = "2016-03-21 08:56:49"
$ws3.Range("F3").Value = "227dfc44-15b6-42d6-a696-ae27c012d273.md"
$ws3.Range("G3").Value = "227dfc44-15b6-42d6-a696-ae27c012d273.49a318ecf0566613ef68ad4b941dea8039fc0f09.de-de.xlf"
$ws3.Range("H3").Value = "2016-03-21 08:57:19"
$ws3.Range("J3").Value = "Include"
$ws3.Range("A4").Value = "3343e4d7-024a-4224-a054-9d6db3bbcd44.md"
$ws3.Range("B4").Value = ".md"
$ws3.Range("C4").Value = "Ready for handoff"
$ws3.Range("D4").Value = "3343e4d7-024a-4224-a054-9d6db3bbcd44.fbdc8a0c47139da47674a6de4d6b421568f325ef.de-de.xlf"
$ws3.Range("E4").Value = "2016-03-21 09:00:11"
$ws3.Range("F4").Value = "3343e4d7-024a-4224-a054-9d6db3bbcd44.md"
$ws3.Range("G4").Value = "3343e4d7-024a-4224-a054-9d6db3bbcd44.fbdc8a0c47139da47674a6de4d6b421568f325ef.de-de.xlf"
$ws3.Range("H4").Value = "2016-03-21 08:59:36"
$ws3.Range("J4").Value = "Include"

# Hyperlinks: clear existing, re-add with original targets + refreshed display text
$ws3.Cells.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/98d1c6365ca3582fab6e8cfc607746a763847726/e2e/3343e4d7-024a-4224-a054-9d6db3bbcd44.md", "", "", "fffff416882d-0faf-47cd-8774-c7dfb82d5749.md")
$ws3.Hyperlinks.Add($ws3.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6bc9fbce954267ab7771c580955286d3fc9f0fc4/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/3343e4d7-024a-4224-a054-9d6db3bbcd44.fbdc8a0c47139da47674a6de4d6b421568f325ef.de-de.xlf", "", "", "227dfc44-15b6-42d6-a696-ae27c012d273.49a318ecf0566613ef68ad4b941dea8039fc0f09.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/2746a788a014fd8b894c271ac81161717acd1d9f/e2e/3343e4d7-024a-4224-a054-9d6db3bbcd44.md", "", "", "227dfc44-15b6-42d6-a696-ae27c012d273.md")
$ws3.Hyperlinks.Add($ws3.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/2d008efee8f1652de645de925be676c84b817445/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/3343e4d7-024a-4224-a054-9d6db3bbcd44.fbdc8a0c47139da47674a6de4d6b421568f325ef.de-de.xlf", "", "", "227dfc44-15b6-42d6-a696-ae27c012d273.49a318ecf0566613ef68ad4b941dea8039fc0f09.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/98d1c6365ca3582fab6e8cfc607746a763847726/e2e/fffff416882d-0faf-47cd-8774-c7dfb82d5749.md", "", "", "ffffff4298431a-edd8-4e9e-8464-d1e2453c7eec.md")
$ws3.Hyperlinks.Add($ws3.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6a100a0e16111973fd8f31dab16272974c9c453e/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/227dfc44-15b6-42d6-a696-ae27c012d273.49a318ecf0566613ef68ad4b941dea8039fc0f09.de-de.xlf", "", "", "227dfc44-15b6-42d6-a696-ae27c012d273.49a318ecf0566613ef68ad4b941dea8039fc0f09.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("F3"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/fa0771df205d0ae5d824129ff4cfcdcc7e6b30d8/e2e/227dfc44-15b6-42d6-a696-ae27c012d273.md", "", "", "227dfc44-15b6-42d6-a696-ae27c012d273.md")
$ws3.Hyperlinks.Add($ws3.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/63efcaaf9872870bd663e592531ab423482db6ea/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/227dfc44-15b6-42d6-a696-ae27c012d273.49a318ecf0566613ef68ad4b941dea8039fc0f09.de-de.xlf", "", "", "227dfc44-15b6-42d6-a696-ae27c012d273.49a318ecf0566613ef68ad4b941dea8039fc0f09.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/98d1c6365ca3582fab6e8cfc607746a763847726/e2e/ffffff4298431a-edd8-4e9e-8464-d1e2453c7eec.md", "", "", "3343e4d7-024a-4224-a054-9d6db3bbcd44.md")
$ws3.Hyperlinks.Add($ws3.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6a100a0e16111973fd8f31dab16272974c9c453e/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/227dfc44-15b6-42d6-a696-ae27c012d273.49a318ecf0566613ef68ad4b941dea8039fc0f09.de-de.xlf", "", "", "3343e4d7-024a-4224-a054-9d6db3bbcd44.fbdc8a0c47139da47674a6de4d6b421568f325ef.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("F4"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/fa0771df205d0ae5d824129ff4cfcdcc7e6b30d8/e2e/227dfc44-15b6-42d6-a696-ae27c012d273.md", "", "", "3343e4d7-024a-4224-a054-9d6db3bbcd44.md")
$ws3.Hyperlinks.Add($ws3.Range("G4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/63efcaaf9872870bd663e592531ab423482db6ea/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/227dfc44-15b6-42d6-a696-ae27c012d273.49a318ecf0566613ef68ad4b941dea8039fc0f09.de-de.xlf", "", "", "3343e4d7-024a-4224-a054-9d6db3bbcd44.fbdc8a0c47139da47674a6de4d6b421568f325ef.de-de.xlf")
